$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.423900127410889
$ws.Range("B1").Value = 2.674579381942749
$ws.Range("C1").Value = 1.828469276428223
$ws.Range("D1").Value = 1.631266236305237
$ws.Range("E1").Value = 1.448257923126221
